$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) from "GossA-HW35.xpc" to "GossA"
$ws.Name = "GossA"

# Append a new data row (row 16) following the same pattern as the existing rows.
$row = 16

# Copy the formatting of the cell above (bold/bordered/centered style) onto A16.
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = 14
$ws.Cells.Item($row, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($row, 3).Value  = 1.012687467361678
$ws.Cells.Item($row, 4).Value  = 0.9733256173515474
$ws.Cells.Item($row, 5).Value  = 1.003347938170875
$ws.Cells.Item($row, 6).Value  = 0.9907697088879488
$ws.Cells.Item($row, 7).Value  = 1.012687467361678
$ws.Cells.Item($row, 8).Value  = 0.9733256173515474
$ws.Cells.Item($row, 9).Value  = 1.002120725760857
$ws.Cells.Item($row, 10).Value = 0.9906243530516571
$ws.Cells.Item($row, 11).Value = 1.001125728766453
$ws.Cells.Item($row, 12).Value = 0.9802829572373952
$ws.Cells.Item($row, 13).Value = 1.012687467361678
$ws.Cells.Item($row, 14).Value = 0.9883367777612111
$ws.Cells.Item($row, 15).Value = 0.9950326829430123
$ws.Cells.Item($row, 16).Value = 0.9942855620735516
